# Add a RESOURCELINK data source row to the "Database" sheet, shifting the
# existing rows (previously starting at row 3) down by one, and move the
# active selection to F3 (matches the author's saved selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# Insert a new row above the current row 3 (ESBCONF), pushing all
# subsequent rows down by one - this mirrors Excel's own "Insert Sheet Rows".
$ws.Rows.Item(3).Insert()

# Give the new row 3 the same "Good" cell style used by its neighbours
# (style index 22 in the original workbook == the built-in "Good" style).
$ws.Rows.Item(3).Style = "Good"

# Populate the new row with the RESOURCELINK connection details.
$ws.Range("A3").Value = "RESOURCELINK"
$ws.Range("B3").Value = "CMTEST"
$ws.Range("C3").Value = "jm08_cmt"
$ws.Range("D3").Formula = '=CONCATENATE( "mqsisetdbparms ",ConfigData!$D$4," -n ",A3," -u ",B3," -p ",C3)'
$ws.Range("F3").Formula = '=CONCATENATE( "mqsicvp ",ConfigData!$D$4," -n ",A3)'
$ws.Range("H3").Value = "Attempt to standardise the DSN across environments."

# The row Insert/Style operations also stamp E3/G3/I3 with the row style but
# leave them empty - clear them fully so the row only carries the cells it
# actually needs (A,B,C,D,F,H), matching the rest of the sheet's rows.
$ws.Range("E3").Clear()
$ws.Range("G3").Clear()
$ws.Range("I3").Clear()

# Move the selection to F3, as recorded in the saved workbook.
$ws.Activate()
$ws.Range("F3").Select()
